$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 234, pushing the existing rows 234:263 down to 237:266.
$ws.Rows("234:236").Insert()

# Constant values shared across this block of rows (same market/category/etc.).
$mercadoId = 5
$mercado = "Macroferia Regional de Talca"
$region = "Maule"
$codreg = 7
$categoriaId = 100112028
$categoria = "Sandia"
$variedad = "Sin especificar"
$unidadComercial = "`$/unidad"
$origen = "Región del Maule"
$kgUnidades = 1
$clasificacion = "Hortaliza"
$fecha = 44946

$rows = @(234, 235, 236)
$calidad = @{234 = "Extra"; 235 = "Primera"; 236 = "Segunda"}
$volumen = @{234 = 1600; 235 = 2500; 236 = 3000}
$precio  = @{234 = 3000; 235 = 2000; 236 = 1500}

foreach ($r in $rows) {
    $ws.Range("A$r").Value = $mercadoId
    $ws.Range("B$r").Value = $mercado
    $ws.Range("C$r").Value = $region
    $ws.Range("D$r").Value = $fecha
    $ws.Range("E$r").Value = $codreg
    $ws.Range("F$r").Value = $categoriaId
    $ws.Range("G$r").Value = $categoria
    $ws.Range("H$r").Value = $variedad
    $ws.Range("I$r").Value = $calidad[$r]
    $ws.Range("J$r").Value = $volumen[$r]
    $ws.Range("K$r").Value = $precio[$r]
    $ws.Range("L$r").Value = $precio[$r]
    $ws.Range("M$r").Value = $precio[$r]
    $ws.Range("N$r").Value = $unidadComercial
    $ws.Range("O$r").Value = $origen
    $ws.Range("P$r").Value = $precio[$r]
    $ws.Range("Q$r").Value = $kgUnidades
    $ws.Range("R$r").Value = $clasificacion
}
